$wb = $excel.ActiveWorkbook

# --- Fats sheet: remove "Turbino" (it's really a raw/turbinado sugar, not a fat) ---
$fats = $wb.Worksheets.Item("Fats")
[void]$fats.Rows.Item(3).Delete()
[void]$fats.Rows.Item(3).Select()

# --- Sugars sheet: add "Raw Sugar" as the new second ingredient ---
$sugars = $wb.Worksheets.Item("Sugars")
[void]$sugars.Rows.Item(2).Insert()
$sugars.Cells.Item(2, 1).Value = "Raw Sugar"

# Sugars ends up as the active sheet/tab, with the new row selected
[void]$sugars.Activate()
[void]$sugars.Rows.Item(2).Select()
